$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 217, shifting existing rows 217-230 down to 218-231
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new weekly data point
$ws.Cells.Item(217, 1).Value = 9
$ws.Cells.Item(217, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(217, 3).Value = "Metropolitana"
$ws.Cells.Item(217, 4).Value = 44746
$ws.Cells.Item(217, 5).Value = 13
$ws.Cells.Item(217, 6).Value = 100112026
$ws.Cells.Item(217, 7).Value = "Haba"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "Primera"
$ws.Cells.Item(217, 10).Value = 43
$ws.Cells.Item(217, 11).Value = 18000
$ws.Cells.Item(217, 12).Value = 18000
$ws.Cells.Item(217, 13).Value = 18000
$ws.Cells.Item(217, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(217, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(217, 16).Value = 720
$ws.Cells.Item(217, 17).Value = 25
$ws.Cells.Item(217, 18).Value = "Hortaliza"

# Match the date-number format used by the rest of column D
$ws.Cells.Item(217, 4).NumberFormat = $ws.Cells.Item(218, 4).NumberFormat
